# 09/12/2024 - preparing for new version
#
# The "Data" sheet contained 61 data rows; rows 28-34 (7 rows) of raw
# measurements were removed, shifting all subsequent rows up by 7 so the
# sheet now ends at row 54 (dimension A1:W54 instead of A1:W61).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Delete rows 28 through 34 (inclusive), shifting the rows below them up.
$ws.Rows("28:34").Delete()
